$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 4).Value = "64.951.53"
$ws.Cells.Item(2, 5).Value = "  +0.56%  "
$ws.Cells.Item(3, 4).Value = "3.441.41"
$ws.Cells.Item(3, 5).Value = "  -1.05%  "
$ws.Cells.Item(4, 5).Value = "  +0.02%  "
$ws.Cells.Item(5, 4).Value = "'574.69"
$ws.Cells.Item(5, 5).Value = "  -1.50%  "
$ws.Cells.Item(6, 4).Value = "'160.67"
$ws.Cells.Item(6, 5).Value = "  +1.39%  "
$ws.Cells.Item(7, 4).Value = "'0.999"
$ws.Cells.Item(7, 5).Value = "  -0.01%  "
$ws.Cells.Item(8, 4).Value = "3.447.63"
$ws.Cells.Item(8, 5).Value = "  -0.96%  "
$ws.Cells.Item(9, 5).Value = "  +7.60%  "
$ws.Cells.Item(10, 4).Value = "'7.29"
$ws.Cells.Item(10, 5).Value = "  -4.28%  "
$ws.Cells.Item(11, 5).Value = "  +0.17%  "
$ws.Cells.Item(12, 4).Value = "'0.438"
$ws.Cells.Item(12, 5).Value = "  -0.85%  "
$ws.Cells.Item(13, 4).Value = "4.031.92"
$ws.Cells.Item(13, 5).Value = "  -1.02%  "
$ws.Cells.Item(14, 5).Value = "  -2.25%  "
$ws.Cells.Item(15, 5).Value = "  +1.72%  "
$ws.Cells.Item(16, 4).Value = "'27.91"
$ws.Cells.Item(16, 5).Value = "  +0.80%  "
$ws.Cells.Item(17, 4).Value = "64.879.95"
$ws.Cells.Item(17, 5).Value = "  +0.41%  "
$ws.Cells.Item(18, 4).Value = "3.376.63"
$ws.Cells.Item(18, 5).Value = "  -2.81%  "
$ws.Cells.Item(19, 4).Value = "'6.34"
$ws.Cells.Item(19, 5).Value = "  -1.68%  "
$ws.Cells.Item(20, 5).Value = "  -1.34%  "
$ws.Cells.Item(21, 4).Value = "'386.02"
$ws.Cells.Item(21, 5).Value = "  -2.84%  "
$ws.Cells.Item(22, 4).Value = "'8.15"
$ws.Cells.Item(22, 5).Value = "  -4.53%  "
$ws.Cells.Item(23, 4).Value = "'72.91"
$ws.Cells.Item(23, 5).Value = "  +0.98%  "
$ws.Cells.Item(24, 4).Value = "'0.542"
$ws.Cells.Item(24, 5).Value = "  -0.64%  "
$ws.Cells.Item(25, 4).Value = "'1.00"
$ws.Cells.Item(25, 5).Value = "  -0.08%  "
$ws.Cells.Item(26, 5).Value = "  +10.31%  "
$ws.Cells.Item(27, 4).Value = "'9.57"
$ws.Cells.Item(27, 5).Value = "  +0.43%  "
$ws.Cells.Item(28, 4).Value = "'0.179"
$ws.Cells.Item(28, 5).Value = "  -1.67%  "
$ws.Cells.Item(29, 5).Value = "  +0.00%  "
$ws.Cells.Item(30, 4).Value = "'6.19"
$ws.Cells.Item(30, 5).Value = "  +5.70%  "
$ws.Cells.Item(31, 4).Value = "'1.43"
$ws.Cells.Item(31, 5).Value = "  +0.45%  "
$ws.Cells.Item(32, 4).Value = "'2.04"
$ws.Cells.Item(32, 5).Value = "  -1.04%  "
$ws.Cells.Item(33, 2).Value = "EthereumClassic"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(33, 4).Value = "'23.59"
$ws.Cells.Item(33, 5).Value = "  -1.13%  "
$ws.Cells.Item(34, 2).Value = "RenderToken"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(34, 4).Value = "'6.51"
$ws.Cells.Item(34, 5).Value = "  -3.29%  "
$ws.Cells.Item(35, 5).Value = "  +0.19%  "
$ws.Cells.Item(36, 4).Value = "'7.05"
$ws.Cells.Item(36, 5).Value = "  +1.42%  "
$ws.Cells.Item(37, 2).Value = "ImmutableX"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(37, 4).Value = "'1.49"
$ws.Cells.Item(37, 5).Value = "  -1.13%  "
$ws.Cells.Item(38, 2).Value = "Monero"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(38, 4).Value = "'161.67"
$ws.Cells.Item(38, 5).Value = "  +1.75%  "
$ws.Cells.Item(39, 2).Value = "Stacks"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(39, 4).Value = "'1.92"
$ws.Cells.Item(39, 5).Value = "  +1.34%  "
$ws.Cells.Item(40, 2).Value = "Maker"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(40, 4).Value = "3.041.45"
$ws.Cells.Item(40, 5).Value = "  +5.29%  "
$ws.Cells.Item(41, 4).Value = "'0.0764"
$ws.Cells.Item(41, 5).Value = "  -2.81%  "
$ws.Cells.Item(42, 4).Value = "'27.20"
$ws.Cells.Item(42, 5).Value = "  -4.71%  "
$ws.Cells.Item(43, 4).Value = "'4.52"
$ws.Cells.Item(43, 5).Value = "  +1.97%  "
$ws.Cells.Item(44, 4).Value = "'42.84"
$ws.Cells.Item(44, 5).Value = "  +1.50%  "
$ws.Cells.Item(45, 5).Value = "  -2.42%  "
$ws.Cells.Item(46, 4).Value = "'0.770"
$ws.Cells.Item(46, 5).Value = "  -1.95%  "
$ws.Cells.Item(47, 4).Value = "'24.62"
$ws.Cells.Item(47, 5).Value = "  +7.84%  "
$ws.Cells.Item(48, 5).Value = "  -2.86%  "
$ws.Cells.Item(49, 4).Value = "'0.868"
$ws.Cells.Item(49, 5).Value = "  +3.70%  "
$ws.Cells.Item(50, 2).Value = "dogwifhat"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(50, 4).Value = "'2.18"
$ws.Cells.Item(50, 5).Value = "  +2.00%  "
$ws.Cells.Item(51, 2).Value = "Cosmos"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(51, 4).Value = "'6.60"
$ws.Cells.Item(51, 5).Value = "  +2.67%  "
